# "separate dept from affiliations"
#
# - "PI hours" gains a new "app" column (F) holding the full affiliation
#   list that used to live in the "dept" column; "dept" (E) now holds just
#   the PI's primary/home department (the first entry of that list).
# - "dept hours" is renamed to "department hours" and its data is replaced
#   with hours/percentages grouped by each PI's primary department.
# - A new sheet "unit(accumulative) hours" is appended, holding the old
#   "dept hours" data (grouped/accumulated across every department listed
#   in a PI's affiliations).

$wb = $excel.ActiveWorkbook
$piSheet = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. "PI hours": split the old combined "dept" column into a primary
#    "dept" column and a new "app" (affiliations) column.
# ---------------------------------------------------------------------

$appValues = @(
    "['ECE', 'CSL']",
    "['PHYS', 'ECE']",
    "['ME', 'AE', 'CSL']",
    "['ABE', 'CSL']",
    "['ECE', 'CSL']",
    "['CEE', 'CSL']",
    "['ECE', 'CSL']"
)
$deptValues = @('ECE', 'PHYS', 'ME', 'ABE', 'ECE', 'CEE', 'ECE')

# New header cell F1 ("app") — clone the formatting of the existing
# header cells (bold, bordered, centered) before setting its text.
$piSheet.Cells.Item(1, 5).Copy($piSheet.Cells.Item(1, 6))
$piSheet.Cells.Item(1, 6).Value = "app"

for ($i = 0; $i -lt 7; $i++) {
    $row = $i + 2
    $piSheet.Cells.Item($row, 5).Value = $deptValues[$i]
    $piSheet.Cells.Item($row, 6).Value = $appValues[$i]
}

# ---------------------------------------------------------------------
# 2. Rename "dept hours" -> "department hours" and replace its data with
#    hours/percentage grouped by each PI's primary department.
# ---------------------------------------------------------------------

$deptSheet = $wb.Worksheets.Item(2)
$deptSheet.Name = "department hours"

# Clear out the old "dept hours" rows (5 rows remain vs. the old 7).
$deptSheet.Rows.Item(8).EntireRow.Clear()
$deptSheet.Rows.Item(7).EntireRow.Clear()

$deptRows = @(
    @('ECE', 55.5, 35.12658227848101),
    @('PHYS', 41, 25.94936708860759),
    @('ME', 35, 22.15189873417722),
    @('ABE', 25.5, 16.13924050632911),
    @('CEE', 1, 0.6329113924050633)
)

for ($i = 0; $i -lt $deptRows.Count; $i++) {
    $row = $i + 2
    $deptSheet.Cells.Item($row, 2).Value = $deptRows[$i][0]
    $deptSheet.Cells.Item($row, 3).Value = $deptRows[$i][1]
    $deptSheet.Cells.Item($row, 4).Value = $deptRows[$i][2]
}

# ---------------------------------------------------------------------
# 3. Add a new "unit(accumulative) hours" sheet at the end, holding the
#    previous "dept hours" data (accumulated across all listed
#    affiliations, not just the primary department).
# ---------------------------------------------------------------------

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$unitSheet = $wb.Worksheets.Add($null, $lastSheet)
$unitSheet.Name = "unit(accumulative) hours"

# Header row — clone formatting from the "PI hours" header row.
$piSheet.Range("B1:D1").Copy($unitSheet.Range("B1:D1"))
$unitSheet.Cells.Item(1, 2).Value = "unit(accumulative)"
$unitSheet.Cells.Item(1, 3).Value = "hours"
$unitSheet.Cells.Item(1, 4).Value = "percentage"

# Row-index column (A) — clone formatting from the "PI hours" sheet.
$piSheet.Range("A2:A8").Copy($unitSheet.Range("A2:A8"))

$unitRows = @(
    @('CSL', 117, 33.33333333333334),
    @('ECE', 96.5, 27.49287749287749),
    @('PHYS', 41, 11.68091168091168),
    @('ME', 35, 9.971509971509972),
    @('AE', 35, 9.971509971509972),
    @('ABE', 25.5, 7.264957264957265),
    @('CEE', 1, 0.2849002849002849)
)

for ($i = 0; $i -lt $unitRows.Count; $i++) {
    $row = $i + 2
    $unitSheet.Cells.Item($row, 1).Value = $i
    $unitSheet.Cells.Item($row, 2).Value = $unitRows[$i][0]
    $unitSheet.Cells.Item($row, 3).Value = $unitRows[$i][1]
    $unitSheet.Cells.Item($row, 4).Value = $unitRows[$i][2]
}

# Restore the originally-active sheet/tab.
$piSheet.Activate()
